# Sertad1-Ar.xlsx update: recompute TPM-based LR communication scores
# and drop "Resolving-Mac" as a valid target (receiving) cluster,
# while keeping it as a sending cluster (rows 11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows that used to represent "Resolving-Mac" as a target
# cluster (old rows 14-17); the sheet shrinks from 17 to 13 rows.
$ws.Rows("14:17").Delete()

# Row 2: ECs -> ECs  becomes  ECs -> ECs
$ws.Cells.Item(2,7).Value = 19.00288166666667
$ws.Cells.Item(2,8).Value = 57.008645
$ws.Cells.Item(2,9).Value = 0.3092280379804411
$ws.Cells.Item(2,10).Value = 0.3092280379804412
$ws.Cells.Item(2,13).Value = 1.268581666666667
$ws.Cells.Item(2,14).Value = 3.805745
$ws.Cells.Item(2,15).Value = 0.06720906924778088
$ws.Cells.Item(2,16).Value = 0.06720906924778086
$ws.Cells.Item(2,17).Value = 24.10670729616945
$ws.Cells.Item(2,18).Value = 216.960365665525
$ws.Cells.Item(2,19).Value = 0.02078292861798288
$ws.Cells.Item(2,20).Value = 0.02078292861798288

# Row 3: ECs -> FAPs  becomes  ECs -> FAPs
$ws.Cells.Item(3,7).Value = 19.00288166666667
$ws.Cells.Item(3,8).Value = 57.008645
$ws.Cells.Item(3,9).Value = 0.3092280379804411
$ws.Cells.Item(3,10).Value = 0.3092280379804412
$ws.Cells.Item(3,15).Value = 0.5468239485455555
$ws.Cells.Item(3,16).Value = 0.5468239485455555
$ws.Cells.Item(3,17).Value = 196.1361021311656
$ws.Cells.Item(3,18).Value = 1765.22491918049
$ws.Cells.Item(3,19).Value = 0.1690932967294598
$ws.Cells.Item(3,20).Value = 0.1690932967294599

# Row 4: ECs -> MuSCs  becomes  ECs -> MuSCs
$ws.Cells.Item(4,7).Value = 19.00288166666667
$ws.Cells.Item(4,8).Value = 57.008645
$ws.Cells.Item(4,9).Value = 0.3092280379804411
$ws.Cells.Item(4,10).Value = 0.3092280379804412
$ws.Cells.Item(4,13).Value = 7.285187000000001
$ws.Cells.Item(4,14).Value = 21.855561
$ws.Cells.Item(4,15).Value = 0.3859669822066636
$ws.Cells.Item(4,16).Value = 0.3859669822066636
$ws.Cells.Item(4,17).Value = 138.4395464805384
$ws.Cells.Item(4,18).Value = 1245.955918324845
$ws.Cells.Item(4,19).Value = 0.1193518126329984
$ws.Cells.Item(4,20).Value = 0.1193518126329984

# Row 5: ECs -> Resolving-Mac  becomes  FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,7).Value = 15.69618166666667
$ws.Cells.Item(5,8).Value = 47.088545
$ws.Cells.Item(5,9).Value = 0.2554191277779661
$ws.Cells.Item(5,10).Value = 0.2554191277779662
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.268581666666667
$ws.Cells.Item(5,14).Value = 3.805745
$ws.Cells.Item(5,15).Value = 0.06720906924778088
$ws.Cells.Item(5,16).Value = 0.06720906924778086
$ws.Cells.Item(5,17).Value = 19.91188829900278
$ws.Cells.Item(5,18).Value = 179.206994691025
$ws.Cells.Item(5,19).Value = 0.01716648184603712
$ws.Cells.Item(5,20).Value = 0.01716648184603712

# Row 6: FAPs -> ECs  becomes  FAPs -> FAPs
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,9).Value = 0.2554191277779661
$ws.Cells.Item(6,10).Value = 0.2554191277779662
$ws.Cells.Item(6,13).Value = 10.32138733333333
$ws.Cells.Item(6,14).Value = 30.964162
$ws.Cells.Item(6,15).Value = 0.5468239485455555
$ws.Cells.Item(6,16).Value = 0.5468239485455555
$ws.Cells.Item(6,17).Value = 162.0063706360322
$ws.Cells.Item(6,18).Value = 1458.05733572429
$ws.Cells.Item(6,19).Value = 0.1396692959856092
$ws.Cells.Item(6,20).Value = 0.1396692959856092

# Row 7: FAPs -> FAPs  becomes  FAPs -> MuSCs
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,9).Value = 0.2554191277779661
$ws.Cells.Item(7,10).Value = 0.2554191277779662
$ws.Cells.Item(7,13).Value = 7.285187000000001
$ws.Cells.Item(7,14).Value = 21.855561
$ws.Cells.Item(7,15).Value = 0.3859669822066636
$ws.Cells.Item(7,16).Value = 0.3859669822066636
$ws.Cells.Item(7,17).Value = 114.3496186276383
$ws.Cells.Item(7,18).Value = 1029.146567648745
$ws.Cells.Item(7,19).Value = 0.09858334994631979
$ws.Cells.Item(7,20).Value = 0.0985833499463198

# Row 8: FAPs -> MuSCs  becomes  MuSCs -> ECs
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,7).Value = 13.073911
$ws.Cells.Item(8,8).Value = 39.221733
$ws.Cells.Item(8,9).Value = 0.2127477252227749
$ws.Cells.Item(8,10).Value = 0.2127477252227749
$ws.Cells.Item(8,13).Value = 1.268581666666667
$ws.Cells.Item(8,14).Value = 3.805745
$ws.Cells.Item(8,15).Value = 0.06720906924778088
$ws.Cells.Item(8,16).Value = 0.06720906924778086
$ws.Cells.Item(8,17).Value = 16.58532380623167
$ws.Cells.Item(8,18).Value = 149.267914256085
$ws.Cells.Item(8,19).Value = 0.01429857659680534
$ws.Cells.Item(8,20).Value = 0.01429857659680534

# Row 9: FAPs -> Resolving-Mac  becomes  MuSCs -> FAPs
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,7).Value = 13.073911
$ws.Cells.Item(9,8).Value = 39.221733
$ws.Cells.Item(9,9).Value = 0.2127477252227749
$ws.Cells.Item(9,10).Value = 0.2127477252227749
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 10.32138733333333
$ws.Cells.Item(9,14).Value = 30.964162
$ws.Cells.Item(9,15).Value = 0.5468239485455555
$ws.Cells.Item(9,16).Value = 0.5468239485455555
$ws.Cells.Item(9,17).Value = 134.9408993925273
$ws.Cells.Item(9,18).Value = 1214.468094532746
$ws.Cells.Item(9,19).Value = 0.1163355511504026
$ws.Cells.Item(9,20).Value = 0.1163355511504027

# Row 10: MuSCs -> ECs  becomes  MuSCs -> MuSCs
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,7).Value = 13.073911
$ws.Cells.Item(10,8).Value = 39.221733
$ws.Cells.Item(10,9).Value = 0.2127477252227749
$ws.Cells.Item(10,10).Value = 0.2127477252227749
$ws.Cells.Item(10,13).Value = 7.285187000000001
$ws.Cells.Item(10,14).Value = 21.855561
$ws.Cells.Item(10,15).Value = 0.3859669822066636
$ws.Cells.Item(10,16).Value = 0.3859669822066636
$ws.Cells.Item(10,17).Value = 95.24588645635701
$ws.Cells.Item(10,18).Value = 857.2129781072131
$ws.Cells.Item(10,19).Value = 0.08211359747556693
$ws.Cells.Item(10,20).Value = 0.08211359747556693

# Row 11: MuSCs -> FAPs  becomes  Resolving-Mac -> ECs
$ws.Cells.Item(11,1).Value = "Resolving-Mac"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,7).Value = 13.67967333333333
$ws.Cells.Item(11,8).Value = 41.03902
$ws.Cells.Item(11,9).Value = 0.2226051090188178
$ws.Cells.Item(11,10).Value = 0.2226051090188178
$ws.Cells.Item(11,13).Value = 1.268581666666667
$ws.Cells.Item(11,14).Value = 3.805745
$ws.Cells.Item(11,15).Value = 0.06720906924778088
$ws.Cells.Item(11,16).Value = 0.06720906924778086
$ws.Cells.Item(11,17).Value = 17.35378279665556
$ws.Cells.Item(11,18).Value = 156.1840451699
$ws.Cells.Item(11,19).Value = 0.01496108218695554
$ws.Cells.Item(11,20).Value = 0.01496108218695553

# Row 12: MuSCs -> MuSCs  becomes  Resolving-Mac -> FAPs
$ws.Cells.Item(12,1).Value = "Resolving-Mac"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,7).Value = 13.67967333333333
$ws.Cells.Item(12,8).Value = 41.03902
$ws.Cells.Item(12,9).Value = 0.2226051090188178
$ws.Cells.Item(12,10).Value = 0.2226051090188178
$ws.Cells.Item(12,13).Value = 10.32138733333333
$ws.Cells.Item(12,14).Value = 30.964162
$ws.Cells.Item(12,15).Value = 0.5468239485455555
$ws.Cells.Item(12,16).Value = 0.5468239485455555
$ws.Cells.Item(12,17).Value = 141.1932070668045
$ws.Cells.Item(12,18).Value = 1270.73886360124
$ws.Cells.Item(12,19).Value = 0.1217258046800838
$ws.Cells.Item(12,20).Value = 0.1217258046800838

# Row 13: MuSCs -> Resolving-Mac  becomes  Resolving-Mac -> MuSCs
$ws.Cells.Item(13,1).Value = "Resolving-Mac"
$ws.Cells.Item(13,4).Value = "MuSCs"
$ws.Cells.Item(13,7).Value = 13.67967333333333
$ws.Cells.Item(13,8).Value = 41.03902
$ws.Cells.Item(13,9).Value = 0.2226051090188178
$ws.Cells.Item(13,10).Value = 0.2226051090188178
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 7.285187000000001
$ws.Cells.Item(13,14).Value = 21.855561
$ws.Cells.Item(13,15).Value = 0.3859669822066636
$ws.Cells.Item(13,16).Value = 0.3859669822066636
$ws.Cells.Item(13,17).Value = 99.65897833224668
$ws.Cells.Item(13,18).Value = 896.93080499022
$ws.Cells.Item(13,19).Value = 0.08591822215177847
$ws.Cells.Item(13,20).Value = 0.08591822215177847

